$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The author solved two new problems ("Find All Anagrams In a String" and
# "Longest Repeating Character Replacement") and recorded them in the
# tracker. In the sheet this shows up as two brand-new "Done" rows being
# inserted right after row 59 (pushing the existing "To Do" rows down by
# three rows), followed by the two new rows being filled in, and a couple of
# extra blank "ID" rows being extended further down.
#
# Because a real Excel "insert rows" leaves stray formatted-but-empty cells
# behind (and because the sheet contains a second, unrelated mini table
# further down that must NOT move), we reproduce the end result directly by
# writing every affected cell (rows 60-71) to its final value instead of
# performing a literal row insert.
# ---------------------------------------------------------------------------

# Save the current ("before") contents of rows 60-65 so we can shift them
# down by three rows (-> rows 63-68) without hard-coding duplicate literals.
$movedRows = @(60, 61, 62, 63, 64, 65)
$saved = @{}
foreach ($r in $movedRows) {
    $saved[$r] = @(
        $ws.Cells.Item($r, 2).Value2,
        $ws.Cells.Item($r, 3).Value2,
        $ws.Cells.Item($r, 4).Value2,
        $ws.Cells.Item($r, 5).Value2
    )
}

# Clear the whole working area (B60:I71) first so no stale values/formatting
# linger in cells that must end up blank (A-only) after the edit.
$ws.Range("B60:I71").Clear()

# Write the two newly solved problems into rows 60 and 61.
$ws.Range("B60").Value = "Sliding Window"
$ws.Range("C60").Value = "Find All Anagrams In a String"
$ws.Range("D60").Value = "Medium"
$ws.Range("E60").Value = "Done"
$ws.Range("F60").Value = 45899
$ws.Range("G60").Value = "O(n)"
$ws.Range("H60").Value = "O(1)"
$ws.Range("I60").Value = "Sliding Window + Hmap"

$ws.Range("B61").Value = "Sliding Window"
$ws.Range("C61").Value = "Longest Repeating Character Replacement"
$ws.Range("D61").Value = "Medium"
$ws.Range("E61").Value = "Done"
$ws.Range("F61").Value = 45899
$ws.Range("G61").Value = "O(n)"
$ws.Range("H61").Value = "O(1)"
$ws.Range("I61").Value = "Sliding Window + Hmap"

# Reuse the existing date-formatted style (column F, row 59) for the two new
# date cells instead of letting Excel synthesize a brand-new number format.
$ws.Range("F59").Copy()
$ws.Range("F60:F61").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rows 62 and 63 stay blank (only the ID column will be (re)filled below).

# Shift the previously-existing "To Do" rows (old rows 61-65) down by three
# rows, landing on rows 64-68.
for ($i = 0; $i -lt $movedRows.Count; $i++) {
    $oldRow = $movedRows[$i]
    if ($oldRow -eq 60) { continue }  # old row 60 was already blank
    $newRow = $oldRow + 3
    $vals = $saved[$oldRow]
    if ($null -ne $vals[0]) { $ws.Cells.Item($newRow, 2).Value = $vals[0] }
    if ($null -ne $vals[1]) { $ws.Cells.Item($newRow, 3).Value = $vals[1] }
    if ($null -ne $vals[2]) { $ws.Cells.Item($newRow, 4).Value = $vals[2] }
    if ($null -ne $vals[3]) { $ws.Cells.Item($newRow, 5).Value = $vals[3] }
}

# Rows 69, 70 and 71 stay blank (only the ID column is filled below) -- these
# extend the existing run of blank "ID only" rows a little further down.

# Re-sequence the ID column (A) for every row in the touched region: the ID
# is simply (row number - 1) as a plain value throughout the sheet.
for ($r = 60; $r -le 71; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# Reflect the author's final cursor position/selection (D61:I61) like the
# saved workbook does.
$ws.Activate()
$ws.Range("D61:I61").Select()
